# Update "想去人数" (interest count) figures in F column on the
# 展览 (Exhibition) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 279
$ws1.Range("F5").Value = 155
$ws1.Range("F6").Value = 81
$ws1.Range("F7").Value = 274
$ws1.Range("F9").Value = 2008
$ws1.Range("F10").Value = 353
$ws1.Range("F11").Value = 4748
$ws1.Range("F12").Value = 87

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 279
$ws4.Range("F7").Value = 155
$ws4.Range("F8").Value = 81
$ws4.Range("F9").Value = 274
$ws4.Range("F13").Value = 2008
$ws4.Range("F14").Value = 353
$ws4.Range("F15").Value = 4749
$ws4.Range("F16").Value = 87
